$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.529953479766846
$ws.Range("B1").Value = 1.752193212509155
$ws.Range("C1").Value = 1.766326069831848
$ws.Range("D1").Value = 2.20918083190918
$ws.Range("E1").Value = 3.332983732223511
